$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal: after "... I will share with a list of relevant journals. "
# append three new runs (each sized 24/24 half-points, i.e. 12pt,
# matching the rest of the paragraph):
#   1) "Also you can look into the series of Journals under MDPI ("
#   2) "https://www.mdpi.com/about/journals"
#   3) ")"
# with no w:rsidR attribute, matching a plain Word paste.
# ------------------------------------------------------------------

$t1 = "Also you can look into the series of Journals under MDPI ("
$t2 = "https://www.mdpi.com/about/journals"
$t3 = ")"

# ---- Step 1: grab a 1-character range that already carries the
# exact formatting we need (sz=24/szCs=24) but no w:rsidR, so that
# anything pasted from it comes out clean. The lone space between
# "Work on the tech reports/papers." and "Look into some relevant
# journals..." has no w:rsidR on it in the source document.
$anchor = $d.Content
$anchor.Find.Execute("Work on the tech reports/papers.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterAnchor = $anchor.End
$fmtSrc = $d.Range($afterAnchor + 1, $afterAnchor + 2)
$fmtSrc.Copy()

# ---- Step 2: stage the new text off in the empty paragraph at the
# very end of the document (isolated from any same-format neighbour)
# so that turning the placeholder into real text via Range.Text does
# not get silently merged into an existing run.
$stageStart = $d.Content.End - 1
$d.Range($stageStart, $stageStart).Paste()

$p1 = $d.Range($stageStart, $stageStart + 1)
$p1.Text = $t1 + $t2 + $t3

$start1 = $stageStart
$start2 = $start1 + $t1.Length
$start3 = $start2 + $t2.Length
$end3   = $start3 + $t3.Length

# ---- Step 3: force the merged staging text back apart into three
# distinct runs by nudging the font size away and back on each
# sub-range (a real property change splits runs; a true no-op does
# not), while keeping the final value identical (12pt) everywhere.
$r1 = $d.Range($start1, $start2)
$r1.Font.Size = 13
$r1.Font.Size = 12

$r2 = $d.Range($start2, $start3)
$r2.Font.Size = 13
$r2.Font.Size = 12

$r3 = $d.Range($start3, $end3)
$r3.Font.Size = 13
$r3.Font.Size = 12

# ---- Step 4: cut the three staged (correctly-formatted, rsid-free)
# runs as one block and paste them where they really belong: right
# after "... I will share with a list of relevant journals. ".
$staged = $d.Range($start1, $end3)
$staged.Cut()

$target = $d.Content
$target.Find.Execute("I will share with a list of relevant journals. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos = $target.End
$d.Range($insPos, $insPos).Paste()
